$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Relocate the hidden "_GoBack" bookmark: it currently sits right
#    after "Benjamín" (inside the backlog table). In the edited
#    document it should sit (collapsed) right after "Con la IA",
#    at the end of the first table, and the text "Benjamín" must
#    become "Benjamin" wrapped with spell-check proofErr markers.
# ------------------------------------------------------------------

# Locate the (unique) "Con la IA" phrase and remember the position
# immediately following it -- that's where the new _GoBack bookmark
# belongs.
$findRng = $d.Content
$null = $findRng.Find.Execute("Con la IA", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $findRng.End

# Insert a temporary placeholder character at that position, bookmark
# the 1-character span (Bookmarks.Add needs a non-collapsed range to
# position correctly), then delete the placeholder character again.
# The bookmark collapses to the correct location and survives.
$placeholder = $d.Range($goBackPos, $goBackPos)
$placeholder.InsertAfter("@")
$placeholderSpan = $d.Range($goBackPos, $goBackPos + 1)
$null = $d.Bookmarks.Add("_GoBack", $placeholderSpan)
$placeholderSpan2 = $d.Range($goBackPos, $goBackPos + 1)
$placeholderSpan2.Text = ""

# ------------------------------------------------------------------
# 2) Find the second "Benjamín" occurrence (the one in the backlog
#    table -- NOT the "Benjamín Miranda Quispe" author byline at the
#    top -- immediately followed by " Miranda") and rewrite that
#    paragraph so the name reads "Benjamin" (no accent), flanked by
#    proofErr spell-check markers, while the old, now stale,
#    "_GoBack" bookmark that used to live right after it is removed
#    as part of the rewrite.
# ------------------------------------------------------------------
$searchRng = $d.Content
$nameStart = -1
$nameEnd = -1
while ($searchRng.Find.Execute("Benjamín", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $afterRng = $d.Range($searchRng.End, $searchRng.End + 15)
    $afterText = $afterRng.Text
    if (-not $afterText.Contains("Quispe")) {
        $nameStart = $searchRng.Start
        $nameEnd = $searchRng.End
        break
    }
    $searchRng.Start = $searchRng.End
    $searchRng.End = $d.Content.End
}

# The target paragraph text is "Benjamín Miranda" -- 8 ("Benjamín")
# + 8 (" Miranda") = 16 characters from $nameStart.
$targetRng = $d.Range($nameStart, $nameStart + 16)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="006F5BA5" w:rsidRDefault="00083043"><w:pPr><w:spacing w:line="256" w:lineRule="auto"/><w:ind w:left="5"/><w:jc w:val="both"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Benjamin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006F5BA5"><w:t xml:space="preserve"> Miranda</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRng.InsertXML($xml)
